# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q3" and "总计", holding
#   the per-fund holding detail for the new quarter (mirrors the layout of
#   "2021-Q3" but with header D1 = "基金规模" instead of "基金金额").
# - Update the "总计" (summary) sheet with a new leading row for 2022-Q1,
#   pushing the existing 2021-Q3 summary row down.

$wb = $excel.ActiveWorkbook

$sheetQ3 = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------
# 1. New sheet "2022-Q1", placed right after "2021-Q3"
# ---------------------------------------------------------------------
$sheetQ1 = $wb.Worksheets.Add($null, $sheetQ3)
$sheetQ1.Name = "2022-Q1"

# NOTE: worksheet handles are positional, so grab "总计" only *after* the
# new sheet has been inserted (insertion shifts its index).
$sheetTotal = $wb.Worksheets.Item("总计")

# Copy the header formatting (font/border/alignment) from "2021-Q3" so the
# new sheet's header row looks the same.
$sheetQ3.Range("B1:H1").Copy($sheetQ1.Range("B1:H1"))

$sheetQ1.Range("B1").Value = "基金代码"
$sheetQ1.Range("C1").Value = "基金名称"
$sheetQ1.Range("D1").Value = "基金规模"
$sheetQ1.Range("E1").Value = "股票总仓位"
$sheetQ1.Range("F1").Value = "仓位占比"
$sheetQ1.Range("G1").Value = "持有市值(亿元)"
$sheetQ1.Range("H1").Value = "仓位排名"

# Row 2: 011243 - A share
$sheetQ3.Range("A2").Copy($sheetQ1.Range("A2"))
$sheetQ1.Range("A2").Value = 0

$sheetQ1.Range("B2").NumberFormat = "@"
$sheetQ1.Range("B2").Value = "011243"
$sheetQ1.Range("C2").Value = "万家惠裕回报6个月持有期混合型证券投资基金A"
$sheetQ1.Range("D2").NumberFormat = "@"
$sheetQ1.Range("D2").Value = "4.93"
$sheetQ1.Range("E2").NumberFormat = "@"
$sheetQ1.Range("E2").Value = "23.04"
$sheetQ1.Range("F2").NumberFormat = "@"
$sheetQ1.Range("F2").Value = "0.83"
$sheetQ1.Range("G2").NumberFormat = "@"
$sheetQ1.Range("G2").Value = "0.0409"
$sheetQ1.Range("H2").Value = 8

# Row 3: 011244 - C share
$sheetQ3.Range("A2").Copy($sheetQ1.Range("A3"))
$sheetQ1.Range("A3").Value = 1

$sheetQ1.Range("B3").NumberFormat = "@"
$sheetQ1.Range("B3").Value = "011244"
$sheetQ1.Range("C3").Value = "万家惠裕回报6个月持有期混合型证券投资基金C"
$sheetQ1.Range("D3").NumberFormat = "@"
$sheetQ1.Range("D3").Value = "0.14"
$sheetQ1.Range("E3").NumberFormat = "@"
$sheetQ1.Range("E3").Value = "23.04"
$sheetQ1.Range("F3").NumberFormat = "@"
$sheetQ1.Range("F3").Value = "0.83"
$sheetQ1.Range("G3").NumberFormat = "@"
$sheetQ1.Range("G3").Value = "0.0012"
$sheetQ1.Range("H3").Value = 8

# ---------------------------------------------------------------------
# 2. Update "总计" sheet: push the old 2021-Q3 row to row 3, and add the
#    new 2022-Q1 row as row 2.
# ---------------------------------------------------------------------
$sheetTotal.Range("A2").Copy($sheetTotal.Range("A3"))
$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("B3").Value = "2021-Q3"
$sheetTotal.Range("C3").Value = 1
$sheetTotal.Range("D3").Value = 0.48

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.04
